$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1)
$ws.Range("A1").Value = "0-0"
$ws.Range("B1").Value = "A-1(1)"
$ws.Range("C1").Value = "A-2(2)"
$ws.Range("D1").Value = "B-1(3)"
$ws.Range("E1").Value = "B-2(4)"
$ws.Range("F1").Value = "B-3(5)"
$ws.Range("G1").Value = "C-1(6)"
$ws.Range("H1").Value = "C-2(7)"
$ws.Range("I1").Value = "C-3(8)"
$ws.Range("J1").Value = "C-4(9)"
$ws.Range("K1").Value = "C-5(10)"

# Row 2
$ws.Range("A2").Value = "A-1(1)"
$ws.Range("B2").Value = "(1-1)"
$ws.Range("C2").Value = "1(1-2)"
$ws.Range("E2").Value = "(1-4)"

# Row 3
$ws.Range("A3").Value = "A-2(2)"

# Row 4
$ws.Range("A4").Value = "B-1(3)"

# Row 5
$ws.Range("A5").Value = "B-2(4)"

# Row 6
$ws.Range("A6").Value = "B-3(5)"

# Row 7
$ws.Range("A7").Value = "C-1(6)"

# Row 8
$ws.Range("A8").Value = "C-2(7)"

# Row 9
$ws.Range("A9").Value = "C-3(8)"

# Row 10
$ws.Range("A10").Value = "C-4(9)"

# Row 11
$ws.Range("A11").Value = "C-5(10)"
$ws.Range("J11").Value = "graph[10][9]"

# Column widths: split former single "1:11" block into 1:9, 10:10, 11:11
# (only column J/10 actually changes width; leave the others untouched so
# their stored width stays exactly 7.625)
$ws.Columns(10).ColumnWidth = 10.08

# Update selection to I6
$ws.Range("I6").Select()
